# DDT-Framework refactor: rename/reshuffle test sheets, add a
# "RemoveCustomerTest" sheet (copy of the trimmed AddCustomerTest),
# trim AddCustomerTest to 2 data rows, lowercase the Y/N run flags,
# and rename test_suite -> testSuite with lowercase tcid/runmode headers.
#
# NOTE: worksheet handles in this engine behave like positional
# references, so any Add/Copy/Delete of a sheet can invalidate handles
# obtained earlier. We therefore perform every edit that only touches
# existing sheets (no sheet insert) FIRST, and only fetch a fresh
# handle for the brand-new sheet right after the single Copy() call,
# which is the last structural change we make.

$wb = $excel.ActiveWorkbook

$addCustomer = $wb.Worksheets.Item("AddCustomerTest")
$openAccount = $wb.Worksheets.Item("OpenAccountTest")
$testSuite   = $wb.Worksheets.Item("test_suite")

# --- 1. Trim AddCustomerTest down to firstname/lastname rows for Joao & Maria ---
# Remove the Jose/Alvez and Jorge/Souza rows entirely.
$addCustomer.Rows.Item(4).EntireRow.Delete()
$addCustomer.Rows.Item(4).EntireRow.Delete()

# Joao's row now flags lower-case "y"; Maria's row flags "n" (no alert).
$addCustomer.Range("E2").Value = "y"
$addCustomer.Range("E3").Value = "n"

# --- 2. OpenAccountTest: drop the wrap-text formatting on the customer name
#        cell, and lower-case its run flag ---
$openAccount.Range("A2").WrapText = $false
$openAccount.Range("C2").Value = "y"

# --- 3. testSuite headers + flags: TCID -> tcid, Runmode -> runmode,
#        Y -> y/n, plus rename the tab ---
$testSuite.Range("A1").Value = "tcid"
$testSuite.Range("B1").Value = "runmode"
$testSuite.Range("B2").Value = "y"
$testSuite.Range("B3").Value = "n"
$testSuite.Range("B4").Value = "y"
$testSuite.Name = "testSuite"

# --- 4. Duplicate the (now trimmed & fixed-up) AddCustomerTest sheet as
#        RemoveCustomerTest, inserted right before it, to match the new
#        first tab. This is the only operation that inserts a sheet, so
#        we do it last and immediately re-resolve the resulting handle. ---
$addCustomer.Copy($addCustomer)
$removeCustomer = $wb.Worksheets.Item("AddCustomerTest (2)")
$removeCustomer.Name = "RemoveCustomerTest"

# --- 5. Leave each sheet's cursor on the cell it was last edited at ---
$removeCustomer.Activate()
$removeCustomer.Range("E7").Select()

$openAccount.Activate()
$openAccount.Range("F9").Select()

$testSuite.Activate()
$testSuite.Range("C11").Select()

$addCustomer.Activate()
$addCustomer.Range("H28").Select()
